$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false
$svc = $wb.Worksheets.Item("individual_services")
$svc.Select()
Write-Host "Before:" $excel.ActiveWindow.ScrollRow
$excel.ActiveWindow.ScrollRow = 5
Write-Host "After:" $excel.ActiveWindow.ScrollRow
